$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously started with a header row (row 1) followed directly by
# data rows. Two new rows are inserted at the very top:
#   - a new row 1 holding a plain numeric column-index sequence (0..13)
#   - a new row 2 that is blank except for "Washer" in column E
# Everything that used to live in rows 1-16 shifts down to rows 3-18.
$ws.Rows("1:2").Insert()

# New row 1: numeric column indexes 0 .. 13 across A1:N1.
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Give the new index row the same look the old header row (now row 3) used to
# have: bold text, thin box border, centered/top-aligned.
$hdr = $ws.Range("A1:N1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# New row 2 is blank apart from "Washer" in column E.
$ws.Range("E2").Value = "Washer"

# The old header row (now row 3) should no longer carry the bold/boxed
# formatting that belonged to row 1 before the insert - only plain text.
$ws.Rows("3:3").ClearFormats()

# The old header row also drops its last two labels ("thread_size" /
# "material_surface" used to live in M3/N3) - those columns are blank now.
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
